# edit.ps1
# Applies two changes described by the target diff:
#
# 1. Slide 6's table (the "SOURCES OF FINANCE" table) gets a new table
#    style GUID: {32C9B805-11E4-4291-8A52-C9B7E005B449} -> {E40F7ED0-61B4-4FB0-83E3-5A112711808C}
#
# 2. The presentation's applied design theme changes from the custom
#    "Integral" theme to the default "Office Theme" (the 10 non-black/white
#    theme colors change; the font scheme and format scheme are identical
#    between the two themes so nothing else needs to move).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{E40F7ED0-61B4-4FB0-83E3-5A112711808C}")
    }
}

# --- 2. Swap the applied theme's colour scheme back to "Office Theme" ----
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (in that RGB() order)
$officeThemeColors = @(
    0,          # dk1      #000000
    16777215,   # lt1      #FFFFFF
    6968388,    # dk2      #44546A
    15132391,   # lt2      #E7E6E6
    13998939,   # accent1  #5B9BD5
    3243501,    # accent2  #ED7D31
    10855845,   # accent3  #A5A5A5
    49407,      # accent4  #FFC000
    12874308,   # accent5  #4472C4
    4697456,    # accent6  #70AD47
    12673797,   # hlink    #0563C1
    7491477     # folHlink #954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeColors[$i - 1]
}
